# Daily attendance processing - 2025-10-12 05:42:12
# Reorder the "Recorded By" (column G) entries so that "System" is listed
# first among the comma-separated recorders, keeping the relative order of
# the remaining entries unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }
    if ($val -notmatch "System") {
        continue
    }

    $parts = $val -split ", "

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Length -eq 0) {
        continue
    }

    $newParts = $systemParts + $otherParts
    $newVal = $newParts -join ", "

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
